$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 247, shifting rows 247:354 down to 248:355
$ws.Rows(247).Insert()

# Populate the newly inserted row 247 with the new data point
$ws.Cells.Item(247, 1).Value = 3
$ws.Cells.Item(247, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(247, 3).Value = "Coquimbo"
$ws.Cells.Item(247, 4).Value = 44726
$ws.Cells.Item(247, 5).Value = 5
$ws.Cells.Item(247, 6).Value = 100112012
$ws.Cells.Item(247, 7).Value = "Espinaca"
$ws.Cells.Item(247, 8).Value = "Sin especificar"
$ws.Cells.Item(247, 9).Value = "Primera"
$ws.Cells.Item(247, 10).Value = 230
$ws.Cells.Item(247, 11).Value = 3500
$ws.Cells.Item(247, 12).Value = 4000
$ws.Cells.Item(247, 13).Value = 3761
$ws.Cells.Item(247, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(247, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(247, 16).Value = 1254
$ws.Cells.Item(247, 17).Value = 3
$ws.Cells.Item(247, 18).Value = "Hortaliza"

# Give D247 the same number format as the date column elsewhere (style index 2)
$ws.Cells.Item(247, 4).NumberFormat = $ws.Cells.Item(248, 4).NumberFormat
